$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Record the "Aktiv Trekk" (active deduction) check: source is the employer
# (AG = Arbeidsgiver) and the withheld amount is 123.
# Shared-string table order matters (new unique strings are appended in the
# order they are first written), so write Y2 before T2 to match the
# upstream commit's string order (AG=41, 123=42).
$ws.Range("Y2").Value = "AG"
$ws.Range("T2").Value = "123"

# Update the view state: scroll the viewport and move the active selection.
$win = $excel.ActiveWindow
$win.ScrollColumn = 10
$win.ScrollRow = 1
$ws.Range("U7").Select()
